$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Row 66 / column B: a second Telegram sticker id was appended to the
#    existing one, joined with " | " (same convention used elsewhere in
#    this sheet for memes that have multiple stickers, e.g. B142).
#    Pull the formatting from B142 (an existing multi-sticker cell) first,
#    then write the final value.
$ws.Range("B142").Copy($ws.Range("B66"))
$ws.Range("B66").Value = "CAACAgEAAxkBAAIWhGAcVHK6lOuJveIQO9PF-bX31isbAAKdAAPFhdoNmGyQjL3bJRoeBA | CAACAgEAAxkBAAIWhmAcVJHpZTL3eU_-TVofgEvepYWiAAKeAAPFhdoNffb48TGnDh8eBA"

# 2) Three new Meme / StickerID rows appended right after the existing
#    data (rows 147-149 already exist as blank, height-only rows).
$newRows = @(
    @{ Row = 147; Meme = "ya ya posi posi | ya posi"; Sticker = "CAACAgEAAxkBAAIk5GA-6aPOAm9NGn_CdKdyUzg02j4jAAIaAAPFhdoNjc3L11hL4mAeBA" },
    @{ Row = 148; Meme = "yo diria que si"; Sticker = "CAACAgEAAxkBAAIk7GA-9wYGYm0sqcrK3O0QS7SurDZ2AALmAgACBTr4RcZstJDezOJ9HgQ" },
    @{ Row = 149; Meme = "esa perra esta loca"; Sticker = "CAACAgEAAxkBAAIk7mA-91HOSLrQZHOyJFtqCt8FnCedAAJeAQAClKP5RbHRpUZ8XrMVHgQ" }
)

foreach ($r in $newRows) {
    # Pull formatting from the row directly above (rows 144-146 use the
    # same "last block" style) so the new rows match.
    $srcRow = $r.Row - 1
    $ws.Range("A" + $srcRow + ":B" + $srcRow).Copy($ws.Range("A" + $r.Row + ":B" + $r.Row))
    $ws.Range("A" + $r.Row).Value = $r.Meme
    $ws.Range("B" + $r.Row).Value = $r.Sticker
}

# 3) Drop the trailing empty formatted row 998 (last row in the sheet).
$ws.Rows.Item(998).Delete()
